# Refactor the AddProduct test-data sheet:
#   - Remove the "PageTitleMyAccountPage" / "My Account" column (old column D).
#     All the columns to its right (ProductCategoryName, ProductSubCategoryName,
#     ProductName) shift one column to the left.
#   - Normalize the formatting on the shifted-in cells so they match the plain
#     bordered style used elsewhere in the table (no stray "apply fill" flag).
#   - Move the active selection to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D ("PageTitleMyAccountPage" header / "My Account" values).
# Excel shifts columns E:G left into D:F automatically.
$ws.Columns("D").Delete()

# Rows 5-8 (old E5:G8, now D5:F8) carried a style with an extraneous
# "applyFill" flag. Reformat them to match the plain style used by the
# rest of the table (as used by, e.g., D2) so the whole table is uniform.
[void]$ws.Range("D2").Copy()
$ws.Range("D5:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection, matching the saved cursor position after the edit.
[void]$ws.Range("C15").Select()
